$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.259.86'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '3.737.65'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''592.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '''166.92'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.66%  '
$ws.Range('D7').Value = '3.736.37'
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '''0.520'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '''0.0000259'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('D14').Value = '''36.17'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.74%  '
$ws.Range('D15').Value = '4.367.66'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '3.748.78'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').Value = '68.227.43'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').Value = '''17.86'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('D19').Value = '''7.00'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').Value = '''10.65'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').Value = '''465.49'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').Value = '''0.696'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').Value = '''0.0000148'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +9.45%  '
$ws.Range('D25').Value = '''83.82'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('D27').Value = '''11.88'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''2.76'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''7.29'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''29.83'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '''2.15'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = '''9.16'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.07%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  --%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.692.84'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.101'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '''3.43'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.138'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').Value = '''0.993'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '''5.78'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = '''44.05'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +17.88%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.300'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '''46.60'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.75%  '
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = '''8.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '''389.75'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''144.16'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.746.54'
$ws.Range('E51').Value = '  +3.27%  '
